$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5725316666666668
$ws.Range("H2").Value = 1.717595
$ws.Range("I2").Value = 0.3864899584549088
$ws.Range("J2").Value = 0.3864899584549088
$ws.Range("M2").Value = 1.716657
$ws.Range("N2").Value = 5.149971
$ws.Range("O2").Value = 0.3840886036988016
$ws.Range("P2").Value = 0.3840886036988015
$ws.Range("Q2").Value = 0.9828404933050001
$ws.Range("R2").Value = 8.845564439745001
$ws.Range("S2").Value = 0.1484463884865538
$ws.Range("T2").Value = 0.1484463884865537
$ws.Range("G3").Value = 0.5725316666666668
$ws.Range("H3").Value = 1.717595
$ws.Range("I3").Value = 0.3864899584549088
$ws.Range("J3").Value = 0.3864899584549088
$ws.Range("O3").Value = 0.07870146593648156
$ws.Range("P3").Value = 0.07870146593648154
$ws.Range("Q3").Value = 0.2013883954377778
$ws.Range("R3").Value = 1.81249555894
$ws.Range("S3").Value = 0.03041732630013118
$ws.Range("T3").Value = 0.03041732630013117
$ws.Range("G4").Value = 0.5725316666666668
$ws.Range("H4").Value = 1.717595
$ws.Range("I4").Value = 0.3864899584549088
$ws.Range("J4").Value = 0.3864899584549088
$ws.Range("M4").Value = 1.677572333333333
$ws.Range("N4").Value = 5.032717
$ws.Range("O4").Value = 0.3753437146230962
$ws.Range("P4").Value = 0.3753437146230962
$ws.Range("Q4").Value = 0.9604632839572224
$ws.Range("R4").Value = 8.644169555615001
$ws.Range("S4").Value = 0.1450665766709916
$ws.Range("T4").Value = 0.1450665766709916
$ws.Range("G5").Value = 0.5725316666666668
$ws.Range("H5").Value = 1.717595
$ws.Range("I5").Value = 0.3864899584549088
$ws.Range("J5").Value = 0.3864899584549088
$ws.Range("M5").Value = 0.7234496666666667
$ws.Range("N5").Value = 2.170349
$ws.Range("O5").Value = 0.1618662157416207
$ws.Range("P5").Value = 0.1618662157416207
$ws.Range("Q5").Value = 0.4141978434061112
$ws.Range("R5").Value = 3.727780590655
$ws.Range("S5").Value = 0.06255966699723231
$ws.Range("T5").Value = 0.0625596669972323
$ws.Range("H6").Value = 0.919331
$ws.Range("I6").Value = 0.2068661122070742
$ws.Range("J6").Value = 0.2068661122070743
$ws.Range("M6").Value = 1.716657
$ws.Range("N6").Value = 5.149971
$ws.Range("O6").Value = 0.3840886036988016
$ws.Range("P6").Value = 0.3840886036988015
$ws.Range("Q6").Value = 0.5260586654889999
$ws.Range("R6").Value = 4.734527989401
$ws.Range("S6").Value = 0.07945491619021476
$ws.Range("T6").Value = 0.07945491619021476
$ws.Range("H7").Value = 0.919331
$ws.Range("I7").Value = 0.2068661122070742
$ws.Range("J7").Value = 0.2068661122070743
$ws.Range("O7").Value = 0.07870146593648156
$ws.Range("P7").Value = 0.07870146593648154
$ws.Range("R7").Value = 0.9701258764120001
$ws.Range("S7").Value = 0.01628066628327743
$ws.Range("T7").Value = 0.01628066628327743
$ws.Range("H8").Value = 0.919331
$ws.Range("I8").Value = 0.2068661122070742
$ws.Range("J8").Value = 0.2068661122070743
$ws.Range("M8").Value = 1.677572333333333
$ws.Range("N8").Value = 5.032717
$ws.Range("O8").Value = 0.3753437146230962
$ws.Range("P8").Value = 0.3753437146230962
$ws.Range("Q8").Value = 0.5140814169252222
$ws.Range("R8").Value = 4.626732752327
$ws.Range("S8").Value = 0.07764589498544147
$ws.Range("T8").Value = 0.07764589498544147
$ws.Range("H9").Value = 0.919331
$ws.Range("I9").Value = 0.2068661122070742
$ws.Range("J9").Value = 0.2068661122070743
$ws.Range("M9").Value = 0.7234496666666667
$ws.Range("N9").Value = 2.170349
$ws.Range("O9").Value = 0.1618662157416207
$ws.Range("P9").Value = 0.1618662157416207
$ws.Range("Q9").Value = 0.2216965685021111
$ws.Range("R9").Value = 1.995269116519
$ws.Range("S9").Value = 0.0334846347481406
$ws.Range("T9").Value = 0.0334846347481406
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5451493333333333
$ws.Range("H10").Value = 1.635448
$ws.Range("I10").Value = 0.3680053968340403
$ws.Range("J10").Value = 0.3680053968340404
$ws.Range("M10").Value = 1.716657
$ws.Range("N10").Value = 5.149971
$ws.Range("O10").Value = 0.3840886036988016
$ws.Range("P10").Value = 0.3840886036988015
$ws.Range("Q10").Value = 0.9358344191119998
$ws.Range("R10").Value = 8.422509772007999
$ws.Range("S10").Value = 0.1413466790236099
$ws.Range("T10").Value = 0.1413466790236099
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5451493333333333
$ws.Range("H11").Value = 1.635448
$ws.Range("I11").Value = 0.3680053968340403
$ws.Range("J11").Value = 0.3680053968340404
$ws.Range("O11").Value = 0.07870146593648156
$ws.Range("P11").Value = 0.07870146593648154
$ws.Range("Q11").Value = 0.1917566414328889
$ws.Range("R11").Value = 1.725809772896
$ws.Range("S11").Value = 0.0289625642033756
$ws.Range("T11").Value = 0.0289625642033756
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5451493333333333
$ws.Range("H12").Value = 1.635448
$ws.Range("I12").Value = 0.3680053968340403
$ws.Range("J12").Value = 0.3680053968340404
$ws.Range("M12").Value = 1.677572333333333
$ws.Range("N12").Value = 5.032717
$ws.Range("O12").Value = 0.3753437146230962
$ws.Range("P12").Value = 0.3753437146230962
$ws.Range("Q12").Value = 0.9145274391351109
$ws.Range("R12").Value = 8.230746952215998
$ws.Range("S12").Value = 0.1381285126490353
$ws.Range("T12").Value = 0.1381285126490353
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5451493333333333
$ws.Range("H13").Value = 1.635448
$ws.Range("I13").Value = 0.3680053968340403
$ws.Range("J13").Value = 0.3680053968340404
$ws.Range("M13").Value = 0.7234496666666667
$ws.Range("N13").Value = 2.170349
$ws.Range("O13").Value = 0.1618662157416207
$ws.Range("P13").Value = 0.1618662157416207
$ws.Range("Q13").Value = 0.3943881034835555
$ws.Range("R13").Value = 3.549492931351999
$ws.Range("S13").Value = 0.05956764095801953
$ws.Range("T13").Value = 0.05956764095801952
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.05723766666666667
$ws.Range("H14").Value = 0.171713
$ws.Range("I14").Value = 0.03863853250397663
$ws.Range("J14").Value = 0.03863853250397663
$ws.Range("M14").Value = 1.716657
$ws.Range("N14").Value = 5.149971
$ws.Range("O14").Value = 0.3840886036988016
$ws.Range("P14").Value = 0.3840886036988015
$ws.Range("Q14").Value = 0.098257441147
$ws.Range("R14").Value = 0.8843169703230001
$ws.Range("S14").Value = 0.01484061999842314
$ws.Range("T14").Value = 0.01484061999842314
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.05723766666666667
$ws.Range("H15").Value = 0.171713
$ws.Range("I15").Value = 0.03863853250397663
$ws.Range("J15").Value = 0.03863853250397663
$ws.Range("O15").Value = 0.07870146593648156
$ws.Range("P15").Value = 0.07870146593648154
$ws.Range("Q15").Value = 0.02013338740844445
$ws.Range("R15").Value = 0.181200486676
$ws.Range("S15").Value = 0.003040909149697352
$ws.Range("T15").Value = 0.003040909149697352
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.05723766666666667
$ws.Range("H16").Value = 0.171713
$ws.Range("I16").Value = 0.03863853250397663
$ws.Range("J16").Value = 0.03863853250397663
$ws.Range("M16").Value = 1.677572333333333
$ws.Range("N16").Value = 5.032717
$ws.Range("O16").Value = 0.3753437146230962
$ws.Range("P16").Value = 0.3753437146230962
$ws.Range("Q16").Value = 0.09602032602455554
$ws.Range("R16").Value = 0.864182934221
$ws.Range("S16").Value = 0.01450273031762783
$ws.Range("T16").Value = 0.01450273031762783
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.05723766666666667
$ws.Range("H17").Value = 0.171713
$ws.Range("I17").Value = 0.03863853250397663
$ws.Range("J17").Value = 0.03863853250397663
$ws.Range("M17").Value = 0.7234496666666667
$ws.Range("N17").Value = 2.170349
$ws.Range("O17").Value = 0.1618662157416207
$ws.Range("P17").Value = 0.1618662157416207
$ws.Range("Q17").Value = 0.04140857087077777
$ws.Range("R17").Value = 0.372677137837
$ws.Range("S17").Value = 0.006254273038228306
$ws.Range("T17").Value = 0.006254273038228306
